$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (municipio-nombre) now uses the refArea / dim / URI-Municipio triple
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column F (aragon) now uses the refArea / dim / URI-Comunidad triple
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "dim"
$ws.Range("F4").Value = "URI-Comunidad"

# The mapping-aragon.xlsx entry in F5 is no longer needed
$ws.Range("F5").Clear()
